{"js": "// Apply the \"1.2.4 -> 1.2.5\" version bump + minor text fixes described in\n// the commit. Each change is a targeted, unambiguous text replacement.\n\n// ---------------------------------------------------------------------\n// 1) Revision-history table: Version / Change / Author / Date row.\n//    These four strings are each unique in the whole document, so a\n//    simple body-wide search is safe.\n// ---------------------------------------------------------------------\nasync function replaceUnique(context, findText, newText) {\n  const results = context.document.body.search(findText, { matchCase: true });\n  results.load(\"items/text\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for '\" + findText + \"', found \" + results.items.length\n    );\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nawait replaceUnique(context, \"1.0\", \"1.2.5\");\nawait replaceUnique(context, \"Creation\", \"Update\");\nawait replaceUnique(context, \"Fabr\u00edcio Ara\u00fajo\", \"Julio Paiva\");\nawait replaceUnique(context, \"09/07/2020\", \"31/05/2023\");\n\n// ---------------------------------------------------------------------\n// 2) Precondition text: fix \"usuario\" -> \"usu\u00e1rio\" and add trailing period.\n// ---------------------------------------------------------------------\nawait replaceUnique(\n  context,\n  \"O usuario devidamente autenticado e na tela inicial de cancelar di\u00e1rias\",\n  \"O usu\u00e1rio devidamente autenticado e na tela inicial de cancelar di\u00e1rias.\"\n);\n\n// ---------------------------------------------------------------------\n// 3) Two \"2. System Exibe a mensagem (MSG102 - Confirmar cancelamento)\"\n//    paragraphs need a period added before the trailing space/suffix.\n//    They are not distinguishable via body.search() alone (one text is a\n//    prefix of the other), so match on exact paragraph text instead.\n// ---------------------------------------------------------------------\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst paraEdits = [\n  {\n    oldText: \"2. System Exibe a mensagem (MSG102 - Confirmar cancelamento) ef[1]\",\n    newText: \"2. System Exibe a mensagem (MSG102 - Confirmar cancelamento). ef[1]\",\n  },\n  {\n    oldText: \"2. System Exibe a mensagem (MSG102 - Confirmar cancelamento) \",\n    newText: \"2. System Exibe a mensagem (MSG102 - Confirmar cancelamento). \",\n  },\n];\n\nfor (const edit of paraEdits) {\n  let matchCount = 0;\n  for (const p of paragraphs.items) {\n    if (p.text === edit.oldText) {\n      p.insertText(edit.newText, Word.InsertLocation.replace);\n      matchCount++;\n    }\n  }\n  if (matchCount !== 1) {\n    throw new Error(\n      \"Expected exactly 1 paragraph match for '\" + edit.oldText + \"', found \" + matchCount\n    );\n  }\n}\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 4) Typo fix inside the EF[1] exception-flow run: \"Solcita\u00e7\u00e3o\" -> \"Solicita\u00e7\u00e3o\".\n// ---------------------------------------------------------------------\nawait replaceUnique(\n  context,\n  \"Solcita\u00e7\u00e3o de di\u00e1ria n\u00e3o pode ser cancelada\",\n  \"Solicita\u00e7\u00e3o de di\u00e1ria n\u00e3o pode ser cancelada\"\n);\n\n// ---------------------------------------------------------------------\n// 5) Remove stray tab character before the closing parenthesis in the\n//    EF[2] exception-flow run (MSG217 text).\n// ---------------------------------------------------------------------\nawait replaceUnique(\n  context,\n  \"solicita\u00e7\u00f5es\\t) para o usu\u00e1rio\",\n  \"solicita\u00e7\u00f5es) para o usu\u00e1rio\"\n);\n", "ps1": "# Apply the \"1.2.4 -> 1.2.5\" version bump + minor text fixes described in\n# the commit, using Word COM interop (Find/Replace + Paragraphs).\n\n$d = $word.ActiveDocument\n\nfunction Replace-UniqueText {\n    param(\n        [string]$FindText,\n        [string]$ReplaceText\n    )\n    $range = $d.Content\n    # wdFindContinue = 1 (Wrap), wdReplaceOne = 1 (Replace)\n    $found = $range.Find.Execute($FindText, $true, $false, $false, $false, $false, $true, 1, $false, $ReplaceText, 1)\n    if (-not $found) {\n        throw \"Find/Replace failed for '$FindText'\"\n    }\n}\n\n# ---------------------------------------------------------------------\n# 1) Revision-history table: Version / Change / Author / Date row.\n#    Each of these strings is unique in the document.\n# ---------------------------------------------------------------------\nReplace-UniqueText \"1.0\" \"1.2.5\"\nReplace-UniqueText \"Creation\" \"Update\"\nReplace-UniqueText \"Fabr\u00edcio Ara\u00fajo\" \"Julio Paiva\"\nReplace-UniqueText \"09/07/2020\" \"31/05/2023\"\n\n# ---------------------------------------------------------------------\n# 2) Precondition text: fix \"usuario\" -> \"usu\u00e1rio\" and add trailing period.\n# ---------------------------------------------------------------------\nReplace-UniqueText \"O usuario devidamente autenticado e na tela inicial de cancelar di\u00e1rias\" \"O usu\u00e1rio devidamente autenticado e na tela inicial de cancelar di\u00e1rias.\"\n\n# ---------------------------------------------------------------------\n# 3) Two \"2. System Exibe a mensagem (MSG102 - Confirmar cancelamento)\"\n#    paragraphs need a period added before the trailing space/suffix.\n#    One text is a prefix of the other, so Find/Replace on the whole\n#    document content could hit the wrong paragraph; match on exact\n#    paragraph text (Range.Text includes the trailing paragraph mark,\n#    so compare against that and trim it off before writing back).\n# ---------------------------------------------------------------------\n$paraEdits = @(\n    @{ Old = \"2. System Exibe a mensagem (MSG102 - Confirmar cancelamento) ef[1]\"; New = \"2. System Exibe a mensagem (MSG102 - Confirmar cancelamento). ef[1]\" },\n    @{ Old = \"2. System Exibe a mensagem (MSG102 - Confirmar cancelamento) \";      New = \"2. System Exibe a mensagem (MSG102 - Confirmar cancelamento). \" }\n)\n\nforeach ($edit in $paraEdits) {\n    $matchCount = 0\n    $paragraphs = $d.Paragraphs\n    for ($i = 1; $i -le $paragraphs.Count; $i++) {\n        $p = $paragraphs.Item($i)\n        if ($p.Range.Text -eq ($edit.Old + \"`r\")) {\n            $r = $p.Range\n            [void]$r.MoveEnd(1, -1) # wdCharacter = 1; exclude the paragraph mark\n            $r.Text = $edit.New\n            $matchCount++\n        }\n    }\n    if ($matchCount -ne 1) {\n        throw \"Expected exactly 1 paragraph match for '$($edit.Old)', found $matchCount\"\n    }\n}\n\n# ---------------------------------------------------------------------\n# 4) Typo fix inside the EF[1] exception-flow run: \"Solcita\u00e7\u00e3o\" -> \"Solicita\u00e7\u00e3o\".\n# ---------------------------------------------------------------------\nReplace-UniqueText \"Solcita\u00e7\u00e3o de di\u00e1ria n\u00e3o pode ser cancelada\" \"Solicita\u00e7\u00e3o de di\u00e1ria n\u00e3o pode ser cancelada\"\n\n# ---------------------------------------------------------------------\n# 5) Remove stray tab character before the closing parenthesis in the\n#    EF[2] exception-flow run (MSG217 text).\n# ---------------------------------------------------------------------\n$tabFind = \"solicita\u00e7\u00f5es\" + [char]9 + \") para o usu\u00e1rio\"\n$tabReplace = \"solicita\u00e7\u00f5es) para o usu\u00e1rio\"\nReplace-UniqueText $tabFind $tabReplace\n\nWrite-Output \"done\"\n"}
